# This script regenerates the "K" column (column G) of the save_data
# sheet so that it reflects strike counts (K) instead of the previous
# "Strike#" derived values. Only column G values change; every other
# column is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new K value
$kValues = @{
    2  = 1
    3  = 2
    5  = 1
    6  = 2
    7  = 1
    8  = 1
    9  = 1
    10 = 2
    11 = 4
    12 = 2
    13 = 0
    14 = 3
    15 = 2
    16 = 1
    17 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 1
    25 = 2
    26 = 1
    27 = 1
    28 = 2
    30 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
